# Auto-generated Excel COM-interop script applying Asura_Profits.xlsx price/profit updates
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# ALC row 5
$ws_ALC.Range("H5").Value = 333.44446
$ws_ALC.Range("I5").Value = 350.125
$ws_ALC.Range("K5").Value = 350.125
$ws_ALC.Range("M5").Value = -235.125

# ALC row 18
$ws_ALC.Range("H18").Value = 769.38464
$ws_ALC.Range("I18").Value = 816.8333
$ws_ALC.Range("J18").Value = 200
$ws_ALC.Range("K18").Value = 816.8333
$ws_ALC.Range("L18").Value = 200
$ws_ALC.Range("M18").Value = -532.8333
$ws_ALC.Range("N18").Value = -768

# ALC row 38
$ws_ALC.Range("H38").Value = 445.86667
$ws_ALC.Range("I38").Value = 74
$ws_ALC.Range("K38").Value = 222
$ws_ALC.Range("M38").Value = 150

# ALC row 81
$ws_ALC.Range("H81").Value = 27000
$ws_ALC.Range("J81").Value = 27000
$ws_ALC.Range("L81").Value = 27000
$ws_ALC.Range("N81").Value = -28996

# ALC row 84
$ws_ALC.Range("H84").Value = 27000
$ws_ALC.Range("J84").Value = 27000
$ws_ALC.Range("L84").Value = 81000
$ws_ALC.Range("N84").Value = -90984

# ALC row 113
$ws_ALC.Range("H113").Value = 3625
$ws_ALC.Range("I113").Value = 2500
$ws_ALC.Range("K113").Value = 2500
$ws_ALC.Range("M113").Value = 754

# ALC row 132
$ws_ALC.Range("H132").Value = 2301.1304
$ws_ALC.Range("I132").Value = 1774.0667
$ws_ALC.Range("J132").Value = 3289.375
$ws_ALC.Range("K132").Value = 5322.2001
$ws_ALC.Range("L132").Value = 9868.125
$ws_ALC.Range("M132").Value = -2792.2001
$ws_ALC.Range("N132").Value = -14928.125

# ALC row 141
$ws_ALC.Range("H141").Value = 5941.8
$ws_ALC.Range("I141").Value = 2617.2693
$ws_ALC.Range("J141").Value = 27551.25
$ws_ALC.Range("K141").Value = 7851.8079
$ws_ALC.Range("L141").Value = 82653.75
$ws_ALC.Range("M141").Value = -2671.8079
$ws_ALC.Range("N141").Value = -93013.75

$ws_ARM = $wb.Worksheets.Item("ARM")
# ARM row 74
$ws_ARM.Range("H74").Value = 882.7143
$ws_ARM.Range("I74").Value = 827.96875
$ws_ARM.Range("J74").Value = 1466.6666
$ws_ARM.Range("K74").Value = 827.96875
$ws_ARM.Range("L74").Value = 1466.6666
$ws_ARM.Range("M74").Value = 46.03125
$ws_ARM.Range("N74").Value = -3214.6666

# ARM row 77
$ws_ARM.Range("H77").Value = 882.7143
$ws_ARM.Range("I77").Value = 827.96875
$ws_ARM.Range("J77").Value = 1466.6666
$ws_ARM.Range("K77").Value = 4139.84375
$ws_ARM.Range("L77").Value = 7333.333000000001
$ws_ARM.Range("M77").Value = 228.15625
$ws_ARM.Range("N77").Value = -16069.333

# ARM row 95
$ws_ARM.Range("H95").Value = 30000
$ws_ARM.Range("J95").Value = 30000
$ws_ARM.Range("L95").Value = 30000
$ws_ARM.Range("N95").Value = -35492

# ARM row 132
$ws_ARM.Range("H132").Value = 2356.889
$ws_ARM.Range("I132").Value = 1944.6666
$ws_ARM.Range("J132").Value = 2686.6667
$ws_ARM.Range("K132").Value = 5833.9998
$ws_ARM.Range("L132").Value = 8060.000100000001
$ws_ARM.Range("M132").Value = -3303.9998
$ws_ARM.Range("N132").Value = -13120.0001

$ws_BSM = $wb.Worksheets.Item("BSM")
# BSM row 35
$ws_BSM.Range("H35").Value = 37000
$ws_BSM.Range("J35").Value = 37000
$ws_BSM.Range("L35").Value = 37000
$ws_BSM.Range("N35").Value = -37620

# BSM row 98
$ws_BSM.Range("H98").Value = 30000
$ws_BSM.Range("J98").Value = 30000
$ws_BSM.Range("L98").Value = 30000
$ws_BSM.Range("N98").Value = -35990

# BSM row 133
$ws_BSM.Range("H133").Value = 57031.125
$ws_BSM.Range("I133").Value = 27709
$ws_BSM.Range("J133").Value = 61220
$ws_BSM.Range("K133").Value = 27709
$ws_BSM.Range("L133").Value = 61220
$ws_BSM.Range("M133").Value = -22649
$ws_BSM.Range("N133").Value = -71340

$ws_CRP = $wb.Worksheets.Item("CRP")
# CRP row 17
$ws_CRP.Range("H17").Value = 0
$ws_CRP.Range("I17").Value = 0
$ws_CRP.Range("J17").Value = 0
$ws_CRP.Range("K17").Value = 0
$ws_CRP.Range("L17").Value = 0
$ws_CRP.Range("M17").Value = $null
$ws_CRP.Range("N17").Value = $null

# CRP row 28
$ws_CRP.Range("H28").Value = 267821.5
$ws_CRP.Range("J28").Value = 267821.5
$ws_CRP.Range("L28").Value = 267821.5
$ws_CRP.Range("N28").Value = -268311.5

# CRP row 31
$ws_CRP.Range("H31").Value = 1717.9344
$ws_CRP.Range("I31").Value = 1254.6
$ws_CRP.Range("J31").Value = 2600.476
$ws_CRP.Range("K31").Value = 1254.6
$ws_CRP.Range("L31").Value = 2600.476
$ws_CRP.Range("M31").Value = -959.5999999999999
$ws_CRP.Range("N31").Value = -3190.476

# CRP row 34
$ws_CRP.Range("H34").Value = 1717.9344
$ws_CRP.Range("I34").Value = 1254.6
$ws_CRP.Range("J34").Value = 2600.476
$ws_CRP.Range("K34").Value = 1254.6
$ws_CRP.Range("L34").Value = 2600.476
$ws_CRP.Range("M34").Value = -1052.6
$ws_CRP.Range("N34").Value = -3004.476

# CRP row 41
$ws_CRP.Range("H41").Value = 8312.4
$ws_CRP.Range("I41").Value = 8312.4
$ws_CRP.Range("K41").Value = 8312.4
$ws_CRP.Range("M41").Value = -7884.4

# CRP row 132
$ws_CRP.Range("H132").Value = 347900.16
$ws_CRP.Range("I132").Value = 437041.62
$ws_CRP.Range("J132").Value = 2477
$ws_CRP.Range("K132").Value = 1311124.86
$ws_CRP.Range("L132").Value = 7431
$ws_CRP.Range("M132").Value = -1308594.86
$ws_CRP.Range("N132").Value = -12491

$ws_CUL = $wb.Worksheets.Item("CUL")
# CUL row 69
$ws_CUL.Range("H69").Value = 562.4
$ws_CUL.Range("I69").Value = 328
$ws_CUL.Range("J69").Value = 1500
$ws_CUL.Range("K69").Value = 984
$ws_CUL.Range("L69").Value = 4500
$ws_CUL.Range("M69").Value = -173
$ws_CUL.Range("N69").Value = -6122

# CUL row 72
$ws_CUL.Range("H72").Value = 562.4
$ws_CUL.Range("I72").Value = 328
$ws_CUL.Range("J72").Value = 1500
$ws_CUL.Range("K72").Value = 2952
$ws_CUL.Range("L72").Value = 13500
$ws_CUL.Range("M72").Value = 1104
$ws_CUL.Range("N72").Value = -21612

# CUL row 109
$ws_CUL.Range("H109").Value = 3296.4167
$ws_CUL.Range("I109").Value = 911.1667
$ws_CUL.Range("J109").Value = 5681.6665
$ws_CUL.Range("K109").Value = 2733.5001
$ws_CUL.Range("L109").Value = 17044.9995
$ws_CUL.Range("M109").Value = -1693.5001
$ws_CUL.Range("N109").Value = -19124.9995

# CUL row 137
$ws_CUL.Range("H137").Value = 25643818
$ws_CUL.Range("I137").Value = 2032.8572
$ws_CUL.Range("J137").Value = 55559236
$ws_CUL.Range("K137").Value = 6098.571599999999
$ws_CUL.Range("L137").Value = 166677708
$ws_CUL.Range("M137").Value = -998.5715999999993
$ws_CUL.Range("N137").Value = -166687908

# CUL row 139
$ws_CUL.Range("H139").Value = 2387.2222
$ws_CUL.Range("I139").Value = 2044.3478
$ws_CUL.Range("K139").Value = 6133.0434
$ws_CUL.Range("M139").Value = -993.0434000000005

$ws_GSM = $wb.Worksheets.Item("GSM")
# GSM row 126
$ws_GSM.Range("H126").Value = 2092.3076
$ws_GSM.Range("I126").Value = 1386.5714
$ws_GSM.Range("J126").Value = 2915.6667
$ws_GSM.Range("K126").Value = 4159.7142
$ws_GSM.Range("L126").Value = 8747.000100000001
$ws_GSM.Range("M126").Value = -1689.7142
$ws_GSM.Range("N126").Value = -13687.0001

# GSM row 132
$ws_GSM.Range("H132").Value = 1920.1538
$ws_GSM.Range("I132").Value = 1311.45
$ws_GSM.Range("K132").Value = 3934.35
$ws_GSM.Range("M132").Value = -1404.35

$ws_LTW = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws_LTW.Range("H7").Value = 5321.1333
$ws_LTW.Range("J7").Value = 5255
$ws_LTW.Range("L7").Value = 5255
$ws_LTW.Range("N7").Value = -5479

# LTW row 126
$ws_LTW.Range("H126").Value = 5321.1333
$ws_LTW.Range("J126").Value = 5255
$ws_LTW.Range("L126").Value = 15765
$ws_LTW.Range("N126").Value = -20705

# LTW row 132
$ws_LTW.Range("H132").Value = 4775.875
$ws_LTW.Range("I132").Value = 4816
$ws_LTW.Range("J132").Value = 4728.4546
$ws_LTW.Range("K132").Value = 14448
$ws_LTW.Range("L132").Value = 14185.3638
$ws_LTW.Range("M132").Value = -11918
$ws_LTW.Range("N132").Value = -19245.3638

$ws_WVR = $wb.Worksheets.Item("WVR")
# WVR row 132
$ws_WVR.Range("H132").Value = 1095.8334
$ws_WVR.Range("J132").Value = 1338.25
$ws_WVR.Range("L132").Value = 4014.75
$ws_WVR.Range("N132").Value = -9074.75
